# quiz updates and announcements
# Duplicate the most recent daily-announcements slide (slide 12) to create
# the next day's slide (Friday, February 14), then update its title and
# body text to match the new announcements.

$p = $ppt.ActivePresentation

# Slide 12 ("Monday, February 10") is the template for the new slide.
$template = $p.Slides.Item(12)
$range = $template.Duplicate()
$s = $range.Item(1)

# --- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Friday, February 14"

# --- Body --------------------------------------------------------------
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$rsquo = [char]0x2019
$hellip = [char]0x2026

$paragraphs = @(
    "Happy Valentines Day?",
    "Quizzes from Wed. will be graded by end of weekend.",
    "Gradescope appears to be working well so far",
    ("I" + $rsquo + "m mad about the due date issue. Any thoughts??"),
    "Interested in a discussion section w/ grad ta?",
    ("Could go over" + $hellip + "solutions / advice for hw problems, proofs, etc."),
    "Remember there is a no laptop policy",
    "Please join the class Piazza asap. Tas are posting common questions there.",
    "Today we will begin continue Greedy Algorithms"
)

$tr.Text = [string]::Join("`r", $paragraphs)

# Indent the two sub-bullets (level 2 = <a:pPr lvl="1"/>)
$tr.Paragraphs(4,1).IndentLevel = 2
$tr.Paragraphs(6,1).IndentLevel = 2

# Paragraph 3: "Gradescope" (bold+italic) + " appears to be working well so far"
$p3 = $tr.Paragraphs(3,1)
$p3.Characters(1,10).Font.Bold = $true
$p3.Characters(1,10).Font.Italic = $true

# Paragraph 6: "Could go over...solutions / advice for " + "hw" + " problems, proofs, etc."
# (no bold/italic formatting on this paragraph besides the indent already applied above)

# Paragraph 7: "Remember there is a " + "no laptop policy" (bold+italic)
$p7 = $tr.Paragraphs(7,1)
$p7.Characters(21,16).Font.Bold = $true
$p7.Characters(21,16).Font.Italic = $true

# Paragraph 8: "Please join the class " + "Piazza" (bold+italic) + " asap. " + "Tas" (unformatted) + " are posting common questions there."
$p8 = $tr.Paragraphs(8,1)
$p8.Characters(23,6).Font.Bold = $true
$p8.Characters(23,6).Font.Italic = $true

# Paragraph 9: "Today we will begin continue " + "Greedy Algorithms" (bold+italic)
$p9 = $tr.Paragraphs(9,1)
$p9.Characters(30,17).Font.Bold = $true
$p9.Characters(30,17).Font.Italic = $true

Write-Output ("New slide count: " + $p.Slides.Count)
